# Add a new "Big_Clock" typography entry (Typography sheet, row 10)
# and its corresponding translation entry (Translation sheet, row 36).
# Commit message: "analog and digital clock"

$wb = $excel.ActiveWorkbook

# --- Typography sheet: new row 10 ---------------------------------------
$wsTypo = $wb.Worksheets.Item("Typography")

# Reset style to "Normal" first so the newly written cells don't inherit
# the column default style index (matches the unstyled cells used by the
# other rows of this table).
$wsTypo.Range("B10:J10").Style = "Normal"

$wsTypo.Range("B10").Value = "Big_Clock"
$wsTypo.Range("C10").Value = "verdana.ttf"
$wsTypo.Range("D10").Value = 100
$wsTypo.Range("E10").Value = 4
$wsTypo.Range("F10").Value = "?"
$wsTypo.Range("H10").Value = "0-9"
$wsTypo.Range("J10").Value = "0123456789 :APM"

# --- Translation sheet: new row 36 ---------------------------------------
$wsTrans = $wb.Worksheets.Item("Translation")

$wsTrans.Range("B36:F36").Style = "Normal"

$wsTrans.Range("B36").Value = "SingleUseId44"
$wsTrans.Range("C36").Value = "Big_Clock"
$wsTrans.Range("D36").Value = "Center"
$wsTrans.Range("E36").Value = "LTR"
$wsTrans.Range("F36").Value = "<>"
